$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.318032698037777
$ws.Range("C2").Value = 0.7852788685634022

$ws.Range("B3").Value = 6.88620145266274
$ws.Range("C3").Value = 0.9016251038884805

$ws.Range("B4").Value = 4.115700080400246
$ws.Range("C4").Value = 0.7967010931424806

$ws.Range("B5").Value = 4.048885414110696
$ws.Range("C5").Value = 0.9973405245089189

$ws.Range("B6").Value = 2.866126328869739
$ws.Range("C6").Value = 0.9687686283208775

$ws.Range("B7").Value = 2.695468174868986
$ws.Range("C7").Value = 0.9984887352256483

$ws.Range("B8").Value = 2.053663236250555
$ws.Range("C8").Value = 0.9976476247153119

$ws.Range("B9").Value = 15.45575288312933
$ws.Range("C9").Value = 0.8152371543542067

$ws.Range("B10").Value = 2.11961902364044
$ws.Range("C10").Value = 0.9936814813423145
